$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.Value = "'248.64"
$r.ClearFormats()
$r = $ws.Range('D3')
$r.Value = "'24.12"
$r.ClearFormats()
$r = $ws.Range('D4')
$r.Value = "'5.959"
$r.ClearFormats()
$r = $ws.Range('D5')
$r.Value = "'0.05881"
$r.ClearFormats()
$r = $ws.Range('D6')
$r.Value = "'3.436"
$r.ClearFormats()
$r = $ws.Range('D7')
$r.Value = "'6.515"
$r.ClearFormats()
$r = $ws.Range('D9')
$r.Value = "'0.7966"
$r.ClearFormats()
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$r = $ws.Range('D10')
$r.Value = "'0.01266"
$r.ClearFormats()
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$r = $ws.Range('D11')
$r.Value = "'0.1473"
$r.ClearFormats()
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$r = $ws.Range('D12')
$r.Value = "'0.07729"
$r.ClearFormats()
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$r = $ws.Range('D13')
$r.Value = "'0.03304"
$r.ClearFormats()
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$r = $ws.Range('D14')
$r.Value = "'0.03012"
$r.ClearFormats()
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$r = $ws.Range('D15')
$r.Value = "'0.09247"
$r.ClearFormats()
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$r = $ws.Range('D16')
$r.Value = "'3.575"
$r.ClearFormats()
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$r = $ws.Range('D17')
$r.Value = "'0.001680"
$r.ClearFormats()
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$r = $ws.Range('D18')
$r.Value = "'0.04781"
$r.ClearFormats()
$ws.Range('E18').Value = '17CoinExTokenCET'
$r = $ws.Range('D19')
$r.Value = "'0.006223"
$r.ClearFormats()
$r = $ws.Range('D20')
$r.Value = "'0.005534"
$r.ClearFormats()
$r = $ws.Range('D21')
$r.Value = "'0.001071"
$r.ClearFormats()
$r = $ws.Range('D22')
$r.Value = "'0.0001503"
$r.ClearFormats()
$r = $ws.Range('D23')
$r.Value = "'3.706"
$r.ClearFormats()
$r = $ws.Range('D25')
$r.Value = "'0.3350"
$r.ClearFormats()
$r = $ws.Range('D27')
$r.Value = "'0.0006282"
$r.ClearFormats()
$r = $ws.Range('D40')
$r.Value = "'0.04376"
$r.ClearFormats()
$r = $ws.Range('D41')
$r.Value = "'0.007044"
$r.ClearFormats()
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$r = $ws.Range('D42')
$r.Value = "'0.003607"
$r.ClearFormats()
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$r = $ws.Range('D43')
$r.Value = "'0.1059"
$r.ClearFormats()
$ws.Range('E43').Value = '42BKEXTokenBKK'
$r = $ws.Range('D44')
$r.Value = "'0.009651"
$r.ClearFormats()
$r = $ws.Range('D45')
$r.Value = "'0.002464"
$r.ClearFormats()
$ws.Range('E45').Value = '44ACDXExchangeACXT'
$r = $ws.Range('D46')
$r.Value = "'0.00005898"
$r.ClearFormats()
$r = $ws.Range('D48')
$r.Value = "'0.9920"
$r.ClearFormats()
$r = $ws.Range('D49')
$r.Value = "'0.1076"
$r.ClearFormats()
$ws.Range('E49').Value = '48BOLOBOLO'
$r = $ws.Range('D50')
$r.Value = "'0.00002104"
$r.ClearFormats()
$r = $ws.Range('D51')
$r.Value = "'0.01012"
$r.ClearFormats()
$ws.Range('E51').Value = '50SpecialPowerGoldSPGWorstin24h'
